$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "moctar"
$ws.Cells.Item(8, 3).Value = "moctar@exemple.com"

$d8 = $ws.Cells.Item(8, 4)
$d8.NumberFormat = "@"
$d8.Value = "1234"
$d8.Style = "Normal"

$ws.Cells.Item(8, 5).Value = "B19"
$ws.Cells.Item(8, 6).Value = "Vice-Président"
